# Auto-generated edit script: update cryptos list values per diff (Sat Feb 10 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.026.16"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").Value = "2.514.69"
$ws.Range("E3").Value = "  +1.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12
$ws.Range("E12").Value = "  +0.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.41%  "

# Row 15
$ws.Range("D15").Value = "2.900.04"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16
$ws.Range("D16").Value = "2.512.71"
$ws.Range("E16").Value = "  +0.76%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("D18").Value = "47.880.37"
$ws.Range("E18").Value = "  +1.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "

# Row 21
$ws.Range("E21").Value = "  +16.95%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0948"
$ws.Range("E22").Value = "  +1.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "

# Row 30
$ws.Range("E30").Value = "  +3.74%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0790"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "

# Row 36
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("E37").Value = "  +0.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.08%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.28%  "

# Row 42
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.87%  "

# Row 45
$ws.Range("D45").Value = "2.003.68"
$ws.Range("E45").Value = "  +1.92%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.68%  "

